# Material.xlsx: normalize the "Directory" column (D) casing.
#   data/XML/  -> data/xml/   (rows 2-17, the "xml import script" rows)
#   Nodegoat/  -> nodegoat/   (rows 18-35, the "csv file for Nodegoat" rows)
# Also: widen column D to fit the new text, and leave the last-used
# selection on E29 (reflecting where the editor ended up after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D17").Value = "data/xml/"
$ws.Range("D18:D35").Value = "nodegoat/"

# Column Format > Width... dialog value that round-trips to the OOXML
# <col width="18"/> recorded in the target workbook.
$ws.Columns("D").ColumnWidth = 17.38

$ws.Range("E29").Select() | Out-Null
